$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "tenant_id" field (comment + model label, with its data-validation /
# select-list template code) is no longer exported/imported, so drop the
# whole column (H) that carries it — header row (comment.tenant_id_lbl,
# column H1) and data row (model.tenant_id_lbl, column H2).
$ws.Range("H1:H2").EntireColumn.Delete()
